$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.996.09'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '2.214.52'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''289.33'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').Value = '''88.01'
$ws.Range('E6').Value = '  +5.04%  '
$ws.Range('D7').Value = '''0.516'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.470'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').Value = '''30.66'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').Value = '  +2.66%  '
$ws.Range('E13').Value = '  +2.68%  '
$ws.Range('D14').Value = '2.557.70'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '2.211.12'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').Value = '39.933.37'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '''11.70'
$ws.Range('E19').Value = '  +12.32%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '''5.80'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '''234.91'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '''2.46'
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('E27').Value = '  +5.04%  '
$ws.Range('D28').Value = '''22.59'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').Value = '''153.36'
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('D31').Value = '''31.99'
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').Value = '''4.96'
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').Value = '''0.0718'
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').Value = '''2.82'
$ws.Range('E36').Value = '  +6.31%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('E38').Value = '  +3.02%  '
$ws.Range('D39').Value = '''15.82'
$ws.Range('E39').Value = '  -0.70%  '
$ws.Range('E40').Value = '  +3.75%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''3.87'
$ws.Range('E41').Value = '  +5.37%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.107.00'
$ws.Range('E42').Value = '  +8.77%  '
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('D45').Value = '''9.95'
$ws.Range('E45').Value = '  +5.83%  '
$ws.Range('D46').Value = '''17.62'
$ws.Range('E46').Value = '  +9.39%  '
$ws.Range('D47').Value = '''2.68'
$ws.Range('E47').Value = '  +3.36%  '
$ws.Range('D48').Value = '2.430.39'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').Value = '''1.46'
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').Value = '''88.60'
$ws.Range('E51').Value = '  -0.26%  '
